$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "BAKSHI MOTORS"
$ws.Range("C21").Value = "Mayapuri"
$ws.Range("A22").Select() | Out-Null
